$d = $word.ActiveDocument

$replacements = @(
    @("56×96=5376", "78×11=858"),
    @("31×86=2666", "11×36=396"),
    @("13×28=364", "69×97=6693"),
    @("94×71=6674", "23×49=1127"),
    @("52×50=2600", "33×34=1122"),
    @("34×66=2244", "87×30=2610"),
    @("67×40=2680", "18×45=810"),
    @("46×90=4140", "46×40=1840"),
    @("95×54=5130", "26×11=286"),
    @("67×37=2479", "68×79=5372"),
    @("89×98=8722", "70×80=5600"),
    @("35×35=1225", "64×13=832"),
    @("99×34=3366", "34×25=850"),
    @("80×80=6400", "66×17=1122"),
    @("69×84=5796", "53×19=1007"),
    @("55×49=2695", "36×64=2304"),
    @("99×59=5841", "65×53=3445"),
    @("99×51=5049", "71×38=2698"),
    @("45×24=1080", "84×30=2520"),
    @("87×40=3480", "25×46=1150"),
    @("82×41=3362", "18×43=774"),
    @("40×50=2000", "32×50=1600"),
    @("64×14=896", "88×89=7832"),
    @("50×11=550", "19×74=1406"),
    @("41×11=451", "19×75=1425")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
